$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75, shifting existing rows 75-122 down to 76-123.
$ws.Rows("75:75").Insert()

# Populate the newly inserted row 75 with the new weekly record.
$ws.Range("A75").Value = 9
$ws.Range("B75").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C75").Value = "Metropolitana"
$ws.Range("D75").Value = 44777
$ws.Range("E75").Value = 13
$ws.Range("F75").Value = 100112022
$ws.Range("G75").Value = "Arveja Verde"
$ws.Range("H75").Value = "Perfection"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 52
$ws.Range("K75").Value = 42000
$ws.Range("L75").Value = 42000
$ws.Range("M75").Value = 42000
$ws.Range("N75").Value = "`$/malla 25 kilos"
$ws.Range("O75").Value = "Provincia de Huasco"
$ws.Range("P75").Value = 1680
$ws.Range("Q75").Value = 25
$ws.Range("R75").Value = "Hortaliza"
